$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.253.11'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.60'
$ws.Range("E3").Value = '  +3.54%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.78'
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5061'
$ws.Range("E7").Value = '  +1.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3925'
$ws.Range("E8").Value = '  +0.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09632'
$ws.Range("E9").Value = '  +1.52%  '

$ws.Range("E10").Value = '  +4.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.90'
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.479'
$ws.Range("E12").Value = '  +1.49%  '

$ws.Range("E13").Value = '  +2.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.883.07'
$ws.Range("E14").Value = '  +3.92%  '

$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.415'
$ws.Range("E16").Value = '  +2.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001131'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.90'
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06606'
$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.65'
$ws.Range("E20").Value = '  +2.96%  '

$ws.Range("E21").Value = '  +0.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.187'
$ws.Range("E22").Value = '  +4.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.311.57'
$ws.Range("E23").Value = '  +1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  +1.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  +3.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.569'
$ws.Range("E26").Value = '  +6.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.089.06'
$ws.Range("E27").Value = '  +3.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.22'
$ws.Range("E28").Value = '  +2.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '158.75'
$ws.Range("E29").Value = '  +1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.60'
$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1066'
$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.065'
$ws.Range("E32").Value = '  +0.96%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.628'
$ws.Range("E33").Value = '  +0.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.625'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06763'
$ws.Range("E35").Value = '  -0.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.500'
$ws.Range("E36").Value = '  +6.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02414'
$ws.Range("E37").Value = '  +4.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2190'
$ws.Range("E38").Value = '  +2.05%  '

$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("E40").Value = '  +1.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6353'
$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("E42").Value = '  +3.51%  '

$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.56'
$ws.Range("E44").Value = '  +4.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5988'
$ws.Range("E45").Value = '  +1.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.274'
$ws.Range("E46").Value = '  -0.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.660'
$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.004'
$ws.Range("E48").Value = '  +2.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.12'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.197'
$ws.Range("E50").Value = '  +1.55%  '

$ws.Range("E51").Value = '  +1.68%  '
